# Iteration v0.7.6b -- Logistic Regression
#
# This script reproduces, via Excel COM automation, a new "0.7.6a"/"0.7.6b"
# iteration pair added to the Results / Steps sheets of Results_summary.xlsx.

$wb = $excel.ActiveWorkbook

$wsResults = $wb.Worksheets.Item(1)   # "Results"
$wsSteps   = $wb.Worksheets.Item(2)   # "Steps"

# ---------------------------------------------------------------------
# Results sheet: the 3rd results block (rows 16-20) gets two more
# iteration columns labelled "0.7.6b" (D) and "0.7.6a" (F), replacing the
# placeholder numeric values that used to sit in the header row, and the
# corresponding Accuracy / FPR / F1 figures are filled in underneath.
# ---------------------------------------------------------------------

$wsResults.Range("D16").Value = "0.7.6b"
$wsResults.Range("F16").Value = "0.7.6a"

$wsResults.Range("D18").Value = 0.973232323232323
$wsResults.Range("F18").Value = 0.922727272727273

$wsResults.Range("D19").Value = 0.0488647581441264
$wsResults.Range("F19").Value = 0.0552813425468904

$wsResults.Range("D20").Value = 0.947860304968028
$wsResults.Range("F20").Value = 0.862162162162162

# ---------------------------------------------------------------------
# Steps sheet: document the two new iterations in rows 26-27, and note
# that 0.7.6b added a Logistic Regression model.
# ---------------------------------------------------------------------

$wsSteps.Cells.Item(26, 1).Value = "0.7.6a"
$wsSteps.Cells.Item(27, 1).Value = "0.7.6b"
$wsSteps.Cells.Item(27, 4).Value = "Logistic Regression"

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping, matching the saved workbook
# state: Steps keeps a selection on D27, Results ends up with F20
# selected and becomes the active (visible) tab.
# ---------------------------------------------------------------------

$wsSteps.Range("D27").Select()
$wsResults.Range("F20").Select()
$wsResults.Activate()
